$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 7 ("Textured surfaces") so the
# new tasks land at rows 7-8 and the existing tasks shift down to rows 9-11.
$ws.Rows("7:8").Insert()

# Fill in the two new tasks.
$ws.Range("A7").Value = "Add full screen support"
$ws.Range("B7").Value = 3

$ws.Range("A8").Value = "Add camera controls to model viewer"
$ws.Range("B8").Value = 5

# Match the author's recorded selection after the edit.
$ws.Range("A32").Select()
